# Contractors.xlsx - "Subs" sheet update
# Adds a new contact (Snyder, David) at R&B Incorporated in row 11, expanding
# the services/emails/phone for that row, and corrects the HVAC trade label
# on row 12 (Allied Service Co / Gable) from "HVAC & plumbing" to "HVAC".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subs")

# --- Row 11: R&B Incorporated -- add new contact + expand details ---
$ws.Range("A11").Value = "Snyder"
$ws.Range("B11").Value = "David"
$ws.Range("D11").Value = "A/C & Heating`nSheet Metal`nUNICO`nParts`nDuct Cleaning"
$ws.Range("E11").Value = "sales@rbincorporated.com`ndave@rbincorporated.com"
$ws.Range("F11").Value = "w: 703-683-1996 `nm: 571-238-9099"

# --- Row 12: Allied Service Co / Gable -- trade label correction ---
$ws.Range("D12").Value = "HVAC"

# --- Refresh the E11 hyperlink so it keeps pointing at the sales@ address
#     but displays the (now multi-line) cell text correctly. The cell's
#     existing hyperlink already targets mailto:sales@rbincorporated.com;
#     re-pointing it keeps the same relationship while recording the
#     friendly display text. ---
$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:sales@rbincorporated.com", "", "", "sales@rbincorporated.com") | Out-Null

# --- View state: frozen-pane top-left cell and the active selection moved
#     up one row (now that row 11 has more content, the sheet was scrolled
#     to show it) ---
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("D12").Select()
